# Refresh the cryptos price/volume snapshot (rows 2-51 of Sheet1) and
# swap the PancakeSwap / Decentraland rows (46 <-> 45 rank), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Price values in column D are plain display strings (e.g. "27.542.03",
# "1.100") rather than real numbers, so each one is entered with a
# leading apostrophe to force Excel to store it as literal text instead
# of silently re-parsing/normalising it as a number (which would, e.g.,
# turn "1.100" into 1.1 and drop the trailing zero).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''27.542.03'
$ws.Cells.Item(2, 5).Value = '  -2.55%  '
$ws.Cells.Item(3, 4).Value = '''1.751.04'
$ws.Cells.Item(3, 5).Value = '  -3.38%  '
$ws.Cells.Item(4, 4).Value = '''1.004'
$ws.Cells.Item(4, 5).Value = '  +0.27%  '
$ws.Cells.Item(5, 4).Value = '''324.56'
$ws.Cells.Item(5, 5).Value = '  -0.05%  '
$ws.Cells.Item(7, 4).Value = '''0.4489'
$ws.Cells.Item(7, 5).Value = '  +3.81%  '
$ws.Cells.Item(8, 4).Value = '''0.3623'
$ws.Cells.Item(8, 5).Value = '  -1.12%  '
$ws.Cells.Item(9, 4).Value = '''0.07488'
$ws.Cells.Item(9, 5).Value = '  -1.86%  '
$ws.Cells.Item(10, 4).Value = '''42.18'
$ws.Cells.Item(10, 5).Value = '  -5.69%  '
$ws.Cells.Item(11, 4).Value = '''1.100'
$ws.Cells.Item(11, 5).Value = '  -3.52%  '
$ws.Cells.Item(12, 5).Value = '  +0.14%  '
$ws.Cells.Item(13, 4).Value = '''20.66'
$ws.Cells.Item(13, 5).Value = '  -5.49%  '
$ws.Cells.Item(14, 4).Value = '''6.039'
$ws.Cells.Item(14, 5).Value = '  -4.07%  '
$ws.Cells.Item(15, 4).Value = '''7.149'
$ws.Cells.Item(15, 5).Value = '  -3.93%  '
$ws.Cells.Item(16, 4).Value = '''1.752.23'
$ws.Cells.Item(16, 5).Value = '  -3.88%  '
$ws.Cells.Item(17, 4).Value = '''92.75'
$ws.Cells.Item(17, 5).Value = '  -1.61%  '
$ws.Cells.Item(18, 4).Value = '''0.00001062'
$ws.Cells.Item(18, 5).Value = '  -1.51%  '
$ws.Cells.Item(19, 4).Value = '''0.06391'
$ws.Cells.Item(19, 5).Value = '  -0.26%  '
$ws.Cells.Item(20, 5).Value = '  +0.11%  '
$ws.Cells.Item(21, 4).Value = '''16.88'
$ws.Cells.Item(21, 5).Value = '  -2.86%  '
$ws.Cells.Item(22, 4).Value = '''5.860'
$ws.Cells.Item(22, 5).Value = '  -5.60%  '
$ws.Cells.Item(23, 4).Value = '''27.588.16'
$ws.Cells.Item(23, 5).Value = '  -2.38%  '
$ws.Cells.Item(24, 4).Value = '''11.21'
$ws.Cells.Item(24, 5).Value = '  -2.99%  '
$ws.Cells.Item(25, 4).Value = '''2.098'
$ws.Cells.Item(25, 5).Value = '  -1.82%  '
$ws.Cells.Item(26, 4).Value = '''161.81'
$ws.Cells.Item(26, 5).Value = '  +0.68%  '
$ws.Cells.Item(27, 4).Value = '''20.42'
$ws.Cells.Item(27, 5).Value = '  -0.87%  '
$ws.Cells.Item(28, 4).Value = '''1.952.69'
$ws.Cells.Item(28, 5).Value = '  -3.64%  '
$ws.Cells.Item(29, 4).Value = '''2.114'
$ws.Cells.Item(29, 5).Value = '  -6.05%  '
$ws.Cells.Item(30, 4).Value = '''125.24'
$ws.Cells.Item(30, 5).Value = '  -3.87%  '
$ws.Cells.Item(31, 4).Value = '''1.081'
$ws.Cells.Item(31, 5).Value = '  -8.55%  '
$ws.Cells.Item(32, 4).Value = '''3.674'
$ws.Cells.Item(32, 5).Value = '  +3.75%  '
$ws.Cells.Item(33, 4).Value = '''0.09023'
$ws.Cells.Item(33, 5).Value = '  -0.88%  '
$ws.Cells.Item(34, 4).Value = '''5.536'
$ws.Cells.Item(34, 5).Value = '  -7.26%  '
$ws.Cells.Item(35, 4).Value = '''11.99'
$ws.Cells.Item(35, 5).Value = '  -6.95%  '
$ws.Cells.Item(36, 4).Value = '''0.02319'
$ws.Cells.Item(36, 5).Value = '  -2.86%  '
$ws.Cells.Item(37, 4).Value = '''0.2089'
$ws.Cells.Item(37, 5).Value = '  -3.11%  '
$ws.Cells.Item(38, 4).Value = '''0.6355'
$ws.Cells.Item(38, 5).Value = '  -2.98%  '
$ws.Cells.Item(39, 5).Value = '  -3.05%  '
$ws.Cells.Item(40, 4).Value = '''4.961'
$ws.Cells.Item(40, 5).Value = '  -4.33%  '
$ws.Cells.Item(41, 4).Value = '''1.209'
$ws.Cells.Item(41, 5).Value = '  +0.61%  '
$ws.Cells.Item(42, 4).Value = '''1.001'
$ws.Cells.Item(42, 5).Value = '  +0.09%  '
$ws.Cells.Item(43, 4).Value = '''1.397'
$ws.Cells.Item(43, 5).Value = '  -2.15%  '
$ws.Cells.Item(44, 4).Value = '''7.773'
$ws.Cells.Item(44, 5).Value = '  -3.23%  '
$ws.Cells.Item(45, 4).Value = '''13.33'
$ws.Cells.Item(45, 5).Value = '  -2.88%  '
$ws.Cells.Item(46, 2).Value = 'PancakeSwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(46, 4).Value = '''3.719'
$ws.Cells.Item(46, 5).Value = '  -0.19%  '
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = '''0.5889'
$ws.Cells.Item(47, 5).Value = '  -2.90%  '
$ws.Cells.Item(48, 4).Value = '''121.31'
$ws.Cells.Item(48, 5).Value = '  -3.71%  '
$ws.Cells.Item(49, 4).Value = '''1.951'
$ws.Cells.Item(49, 5).Value = '  -2.85%  '
$ws.Cells.Item(50, 4).Value = '''1.157'
$ws.Cells.Item(50, 5).Value = '  -1.12%  '
$ws.Cells.Item(51, 4).Value = '''0.06876'
$ws.Cells.Item(51, 5).Value = '  -1.68%  '
